# Apply balance sheet quarterly update: shift quarterly columns left by one
# quarter (drop oldest 1400/06 quarter, add new 1401/09 quarter), update
# copyright year, and resize the "current quarter" highlighted column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / label updates ---
$ws.Range("B3").Value = "Copyright @2015 - 2023"
$ws.Range("D8").Value = "فصل سوم منتهی به 1400/09"
$ws.Range("E8").Value = "فصل چهارم منتهی به 1400/12"
$ws.Range("F8").Value = "فصل اول منتهی به 1401/03"
$ws.Range("G8").Value = "فصل دوم منتهی به 1401/06"
$ws.Range("H8").Value = "فصل سوم منتهی به 1401/09"
$ws.Range("D9").Value = "1400-10-30"
$ws.Range("E9").Value = "1401-11-02 (7)"
$ws.Range("F9").Value = "1401-04-29"
$ws.Range("G9").Value = "1401-09-16 (3)"
$ws.Range("H9").Value = "1401-11-02"
$ws.Range("D12").Value = 475987
$ws.Range("E12").Value = 1868143
$ws.Range("F12").Value = 415393
$ws.Range("G12").Value = 283975
$ws.Range("H12").Value = 678035
$ws.Range("D13").Value = 315189
$ws.Range("E13").Value = 261339
$ws.Range("F13").Value = 307450
$ws.Range("G13").Value = 272224
$ws.Range("H13").Value = 255106
$ws.Range("D14").Value = 9828553
$ws.Range("E14").Value = 11169725
$ws.Range("F14").Value = 9669726
$ws.Range("G14").Value = 12498860
$ws.Range("H14").Value = 14783315
$ws.Range("D15").Value = 4991497
$ws.Range("E15").Value = 4608522
$ws.Range("F15").Value = 5278560
$ws.Range("G15").Value = 5337696
$ws.Range("H15").Value = 6348304
$ws.Range("D16").Value = 1224791
$ws.Range("E16").Value = 1409203
$ws.Range("F16").Value = 1598580
$ws.Range("G16").Value = 1777074
$ws.Range("H16").Value = 2740270
$ws.Range("D18").Value = 16836017
$ws.Range("E18").Value = 19316932
$ws.Range("F18").Value = 17269709
$ws.Range("G18").Value = 20169829
$ws.Range("H18").Value = 24805030
$ws.Range("D20").Value = 528178
$ws.Range("F20").Value = 2406699
$ws.Range("G20").Value = 3678659
$ws.Range("H20").Value = 2406698
$ws.Range("F21").Value = 662825
$ws.Range("G21").Value = 866660
$ws.Range("D22").Value = 2732480
$ws.Range("E22").Value = 3300378
$ws.Range("F22").Value = 3279883
$ws.Range("G22").Value = 3292447
$ws.Range("H22").Value = 3319803
$ws.Range("D23").Value = 28361
$ws.Range("E23").Value = 27944
$ws.Range("F23").Value = 28225
$ws.Range("G23").Value = 27690
$ws.Range("H23").Value = 26700
$ws.Range("D26").Value = 3289019
$ws.Range("E26").Value = 3856500
$ws.Range("F26").Value = 6377632
$ws.Range("G26").Value = 7865456
$ws.Range("H26").Value = 6619861
$ws.Range("D27").Value = 20125036
$ws.Range("E27").Value = 23173432
$ws.Range("F27").Value = 23647341
$ws.Range("G27").Value = 28035285
$ws.Range("H27").Value = 31424891
$ws.Range("D29").Value = 4426796
$ws.Range("E29").Value = 5072864
$ws.Range("F29").Value = 5059166
$ws.Range("G29").Value = 6601920
$ws.Range("H29").Value = 8962503
$ws.Range("D32").Value = 114820
$ws.Range("E32").Value = 102254
$ws.Range("F32").Value = 19694
$ws.Range("G32").Value = 29649
$ws.Range("H32").Value = 128979
$ws.Range("D33").Value = 1991159
$ws.Range("E33").Value = 311296
$ws.Range("F33").Value = 309669
$ws.Range("G33").Value = 2047674
$ws.Range("H33").Value = 1779005
$ws.Range("D34").Value = 7627666
$ws.Range("E34").Value = 8302390
$ws.Range("F34").Value = 8482215
$ws.Range("G34").Value = 8738322
$ws.Range("H34").Value = 9173922
$ws.Range("D35").Value = 198741
$ws.Range("E35").Value = 222012
$ws.Range("F35").Value = 366336
$ws.Range("G35").Value = 406685
$ws.Range("H35").Value = 1029437
$ws.Range("D37").Value = 14359182
$ws.Range("E37").Value = 14010816
$ws.Range("F37").Value = 14237080
$ws.Range("G37").Value = 17824250
$ws.Range("H37").Value = 21073846
$ws.Range("E40").Value = 2361111
$ws.Range("F40").Value = 2319444
$ws.Range("G40").Value = 2324523
$ws.Range("H40").Value = 2275844
$ws.Range("D41").Value = 134304
$ws.Range("E41").Value = 139607
$ws.Range("F41").Value = 224323
$ws.Range("G41").Value = 225629
$ws.Range("H41").Value = 228559
$ws.Range("D42").Value = 134304
$ws.Range("E42").Value = 2500718
$ws.Range("F42").Value = 2543767
$ws.Range("G42").Value = 2550152
$ws.Range("H42").Value = 2504403
$ws.Range("D43").Value = 14493486
$ws.Range("E43").Value = 16511534
$ws.Range("F43").Value = 16780847
$ws.Range("G43").Value = 20374402
$ws.Range("H43").Value = 23578249
$ws.Range("E45").Value = 2000000
$ws.Range("D47").Value = 981288
$ws.Range("E47").Value = 0
$ws.Range("E50").Value = 200000
$ws.Range("D56").Value = 3550262
$ws.Range("E56").Value = 4461898
$ws.Range("F56").Value = 4666494
$ws.Range("G56").Value = 5460883
$ws.Range("H56").Value = 5646642
$ws.Range("D57").Value = 5631550
$ws.Range("E57").Value = 6661898
$ws.Range("F57").Value = 6866494
$ws.Range("G57").Value = 7660883
$ws.Range("H57").Value = 7846642
$ws.Range("D58").Value = 20125036
$ws.Range("E58").Value = 23173432
$ws.Range("F58").Value = 23647341
$ws.Range("G58").Value = 28035285
$ws.Range("H58").Value = 31424891

# --- Column widths ---

# Column width adjustments: the highlighted (wider) column moves from F to E
# as the data shifts left by one quarter.
$ws.Columns.Item(4).ColumnWidth = 28.166666666666668
$ws.Columns.Item(5).ColumnWidth = 30.166666666666668
$ws.Columns.Item(6).ColumnWidth = 28.166666666666668
$ws.Columns.Item(7).ColumnWidth = 28.166666666666668
$ws.Columns.Item(8).ColumnWidth = 28.166666666666668
